$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (E1, F1) matching the existing header style
$ws.Range("E1").Value = "Macro F1"
$ws.Range("F1").Value = "Accuracy"
$ws.Range("B1").Copy()
$ws.Range("E1:F1").PasteSpecial(-4122)  # xlPasteFormats

# Update row 2 metric values
$ws.Range("B2").Value = 0.7589023310417855
$ws.Range("C2").Value = 0.7180442374854482
$ws.Range("D2").Value = 0.7369419147848194
$ws.Range("E2").Value = 0.4989854475662059
$ws.Range("F2").Value = 0.5021814598311271

# Row 3: rename model and add the new metric values
$ws.Range("A3").Value = "XGBoost"
$ws.Range("E3").Value = 0.4887298341363462
$ws.Range("F3").Value = 0.4989529440748953
